$wb = $excel.ActiveWorkbook

# The change happens on the "Repayment schedule" worksheet: a new blank
# column is inserted before column N (14), shifting the old "Late",
# "heading" (Over Due) and "Outstanding" columns one place to the right.
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a whole new column at position 14 (N), pushing N->O, O->P, P->Q.
$ws.Columns.Item(14).Insert()

# The newly inserted column keeps a manual width of 11 (matching the
# width of the "In Advance" column to its left) instead of the bestFit
# auto-width the shifted columns retain.
$ws.Columns.Item(14).ColumnWidth = 10.1666667

# Update the sheet's selection/active cell and make it the active tab,
# which moves the workbook's activeTab away from "Edit Repayment
# Schedule" (previously active) onto this sheet.
$ws.Activate()
$ws.Range("J16").Select()
